$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the CPU threshold table (rows 8-9) to match the Disco table (rows 13-14):
# D8: 0.9 -> "70% - 80%" (text, keeping percentage-style look)
# D9: 1   -> 0.81
$ws.Range("D8").Value = "70% - 80%"
$ws.Range("D9").Value = 0.81

# Update the active cell selection to A11 (mirrors last selection state in the saved file)
$ws.Range("A11").Select()

$wb.Save()
